# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, centered, bordered) from AC1
# onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-49: Wins=76, Losses=86, Ties=0
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 76  # AD
    $ws.Cells.Item($r, 31).Value = 86  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
